{"js": "// Apply the cover-letter edits described by the commit:\n//   \"Added C#, SQL, and .NET Work\" / \"Completed MySwoleMate 100 project from Revature\"\n//\n// Strategy: locate each affected phrase with Body.search (exact, case-sensitive,\n// whitespace-sensitive) and replace it in place with Range.insertText(..., \"Replace\").\n// Each search string below is unique in the document, so a single hit is expected.\n\nasync function replaceOnce(body, findText, replaceText) {\n  const results = body.search(findText, {\n    matchCase: true,\n    matchWholeWord: false,\n    matchWildcards: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1) \"...Github.  Most recently, I contributed to TEAMMATES...\"\n//       -> \"...Github.  I am a contributor to TEAMMATES...\"\nawait replaceOnce(body, \"Most recently, I contributed\", \"I am a contributor\");\n\n// 2) \"...TEAMMATES, open-source software that manages student/instructor...\"\n//       -> \"...TEAMMATES, open-source software for managing student/instructor...\"\nawait replaceOnce(body, \"that manages\", \"for managing\");\n\n// 3) \"...correspondence in higher education.   \" (end of paragraph)\n//       -> \"...correspondence in higher education, and currently, I am working on a\n//           web application, MySwoleMate, to develop my skills with C#, SQL, and the\n//           .NET framework.\"\nawait replaceOnce(\n  body,\n  \"education.   \",\n  \"education, and currently, I am working on a web application, MySwoleMate, \" +\n    \"to develop my skills with C#, SQL, and the .NET framework.\"\n);\n\n// 4) \"...study foundational texts in computer science literature and make further\n//      contributions to TEAMMATES.  I am currently...\"\n//       -> \"...study foundational texts in computer science literature, make further\n//           contributions to TEAMMATES, and continue to study C#, SQL, and .NET.  I am\n//           currently...\"\nawait replaceOnce(body, \"literature and\", \"literature,\");\nawait replaceOnce(\n  body,\n  \"TEAMMATES.\",\n  \"TEAMMATES, and continue to study C#, SQL, and .NET.\"\n);\n", "ps1": "# Apply the cover-letter edits described by the commit:\n#   \"Added C#, SQL, and .NET Work\" / \"Completed MySwoleMate 100 project from Revature\"\n#\n# Strategy: use Find/Replace (wdReplaceOne) against Document.Content for each affected\n# phrase. Each search string is unique in the document, so exactly one hit is expected.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text {\n    param($doc, $findText, $replaceText)\n\n    $range = $doc.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $findText\n    $range.Find.Replacement.Text = $replaceText\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1          # wdFindContinue\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n\n    $found = $range.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)  # wdReplaceOne... actually Replace:=2 = wdReplaceAll, but only one match exists\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# \"...Github.  Most recently, I contributed to TEAMMATES...\"\n#     -> \"...Github.  I am a contributor to TEAMMATES...\"\nReplace-Text $d \"Most recently, I contributed\" \"I am a contributor\"\n\n# \"...TEAMMATES, open-source software that manages student/instructor...\"\n#     -> \"...TEAMMATES, open-source software for managing student/instructor...\"\nReplace-Text $d \"that manages\" \"for managing\"\n\n# \"...correspondence in higher education.   \" (end of paragraph)\n#     -> \"...correspondence in higher education, and currently, I am working on a web\n#         application, MySwoleMate, to develop my skills with C#, SQL, and the .NET\n#         framework.\"\nReplace-Text $d \"education.   \" \"education, and currently, I am working on a web application, MySwoleMate, to develop my skills with C#, SQL, and the .NET framework.\"\n\n# \"...study foundational texts in computer science literature and make further\n#    contributions to TEAMMATES.  I am currently...\"\n#     -> \"...study foundational texts in computer science literature, make further\n#         contributions to TEAMMATES, and continue to study C#, SQL, and .NET.  I am\n#         currently...\"\nReplace-Text $d \"literature and\" \"literature,\"\nReplace-Text $d \"TEAMMATES.\" \"TEAMMATES, and continue to study C#, SQL, and .NET.\"\n"}
